$d = $word.ActiveDocument

# 1) Merge the leading space run with "Nome: <nome>" into a single run.
$d.Content.Find.Execute(" Nome: <nome>", $true, $false, $false, $false, $false, $true, 1, $false, " Nome: <nome>", 2) | Out-Null

# 2) Merge the leading space run with "Matricula: <matricula>" into a single run.
$d.Content.Find.Execute(" Matricula: <matricula>", $true, $false, $false, $false, $false, $true, 1, $false, " Matricula: <matricula>", 2) | Out-Null

# 3) Merge the leading space run with "Turma: <turma>" into a single run.
$d.Content.Find.Execute(" Turma: <turma>", $true, $false, $false, $false, $false, $true, 1, $false, " Turma: <turma>", 2) | Out-Null

# 4) Split "Professor: <professor>" so that "Professor: " moves into the
#    first (previously space-only) run, leaving only "<professor>" in the
#    second run.
$professorParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Professor:*") {
        $professorParaIndex = $i
        break
    }
}

$pProfessor = $d.Paragraphs.Item($professorParaIndex)
$rProfessor = $pProfessor.Range
$profStart = $rProfessor.Start

$rProfessorRest = $d.Range($profStart + 1, $rProfessor.End - 1)
$rProfessorRest.Text = "<professor>"

$rProfessorSpace = $d.Range($profStart, $profStart + 1)
$rProfessorSpace.Text = " Professor: "

# Touch formatting on the first piece (no net visual change) so the engine
# keeps the two pieces as separate runs instead of re-merging them.
$rProfessorBoundary = $d.Range($profStart, $profStart + 12)
$rProfessorBoundary.Font.Bold = $true
$rProfessorBoundary.Font.Bold = $false

# 5) Grow the row right below the "Professor" row (trHeight 230 -> 489 twips),
#    and 6) center the (empty) paragraph in that row's first cell.
$table = $d.Tables.Item(1)
$targetRowIndex = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows.Item($i)
    if ($row.Height -and [Math]::Abs($row.Height - (230 / 20.0)) -lt 0.01) {
        $targetRowIndex = $i
        break
    }
}

$row = $table.Rows.Item($targetRowIndex)
$row.Height = 489 / 20.0

$cell = $row.Cells.Item(1)
$cell.Range.ParagraphFormat.Alignment = 1
